$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") from 45188 to 45189 for data rows 2..344
for ($i = 2; $i -le 344; $i++) {
    $cell = $ws.Cells.Item($i, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}

# Row 344 gains an explicit 15pt custom row height
$ws.Rows.Item(344).RowHeight = 15

# Add new row 345 with the new cleared-case entry
$ws.Range("A345").Value2 = "A 44310-2023"

$ws.Range("B345").Value2 = 45187
$ws.Range("B345").NumberFormat = "YYYY-MM-DD"

$ws.Range("C345").Value2 = 45189
$ws.Range("C345").NumberFormat = "YYYY-MM-DD"

$ws.Range("D345").Value2 = "VÄSTERBOTTENS LÄN"
$ws.Range("E345").Value2 = "NORSJÖ"

$ws.Range("G345").Value2 = 1.9
$ws.Range("H345").Value2 = 0
$ws.Range("I345").Value2 = 0
$ws.Range("J345").Value2 = 0
$ws.Range("K345").Value2 = 0
$ws.Range("L345").Value2 = 0
$ws.Range("M345").Value2 = 0
$ws.Range("N345").Value2 = 0
$ws.Range("O345").Value2 = 0
$ws.Range("P345").Value2 = 0
$ws.Range("Q345").Value2 = 0

$ws.Range("R345").WrapText = $true

Write-Host "Edit complete"
